$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.978.45"
$ws.Range("E2").Value = "  -1.96%  "

$ws.Range("D3").Value = "2.391.48"
$ws.Range("E3").Value = "  -3.50%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.09"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.04"
$ws.Range("E6").Value = "  -4.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.532"
$ws.Range("E7").Value = "  -3.58%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -3.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0842"
$ws.Range("E10").Value = "  -3.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.39"
$ws.Range("E11").Value = "  -5.90%  "

$ws.Range("E12").Value = "  -1.61%  "

$ws.Range("D13").Value = "2.767.62"
$ws.Range("E13").Value = "  -3.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.64"
$ws.Range("E14").Value = "  -4.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.28"
$ws.Range("E15").Value = "  -2.73%  "

$ws.Range("D16").Value = "2.405.71"
$ws.Range("E16").Value = "  -2.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.771"
$ws.Range("E17").Value = "  -3.26%  "

$ws.Range("D18").Value = "40.980.42"
$ws.Range("E18").Value = "  -1.82%  "

$ws.Range("D19").Value = "0.0₃0918"
$ws.Range("E19").Value = "  -3.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  -4.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.50"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.90"
$ws.Range("E22").Value = "  -3.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.75"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.66"
$ws.Range("E24").Value = "  -3.62%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -5.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.03"
$ws.Range("E27").Value = "  -3.27%  "

$ws.Range("E28").Value = "  -1.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.41"
$ws.Range("E29").Value = "  -3.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.26"
$ws.Range("E30").Value = "  -6.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.00"
$ws.Range("E31").Value = "  -2.58%  "

$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.24"
$ws.Range("E33").Value = "  -4.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0738"
$ws.Range("E34").Value = "  -3.66%  "

$ws.Range("E35").Value = "  -4.94%  "

$ws.Range("E36").Value = "  -2.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.83"
$ws.Range("E37").Value = "  -3.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.15"
$ws.Range("E38").Value = "  -7.93%  "

$ws.Range("E39").Value = "  -3.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.74"
$ws.Range("E40").Value = "  -7.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.87"
$ws.Range("E41").Value = "  -3.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("E42").Value = "  -7.75%  "

$ws.Range("D43").Value = "1.976.60"
$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("E44").Value = "  -4.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.79"
$ws.Range("E45").Value = "  -7.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.70"
$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.79"
$ws.Range("E47").Value = "  -7.21%  "

$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.18"
$ws.Range("E49").Value = "  -4.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.10"
$ws.Range("E50").Value = "  -2.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.63"
$ws.Range("E51").Value = "  -1.86%  "
